$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.459.26"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.616.30"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.78"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.77%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "1.844.19"
$ws.Range("D13").Value = "1.626.10"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.549"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.85"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "27.435.57"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.71"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.16"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("E24").Value = "  +5.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.76"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0483"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").Value = "1.466.74"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.01%  "
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.948"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.558"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  -3.16%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.96"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.70%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.984"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.71%  "
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.26"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.77%  "
$ws.Range("D47").Value = "1.755.45"
$ws.Range("E47").Value = "  -1.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.41"
$ws.Range("D49").ClearFormats()
$ws.Range("E50").Value = "  -2.47%  "
$ws.Range("E51").Value = "  +1.43%  "
